$wb = $excel.ActiveWorkbook

# Add the new "Hours" worksheet after the existing "Login" sheet
$loginSheet = $wb.Worksheets.Item("Login")
$ws = $wb.Worksheets.Add([System.Type]::Missing, $loginSheet)
$ws.Name = "Hours"

# Populate row 1 with the hours data
$ws.Range("A1").Value = 8
$ws.Range("B1").Value = 8
$ws.Range("C1").Value = 8
$ws.Range("D1").Value = 8
$ws.Range("E1").Value = 8

# Select cell I3 on the Hours sheet and make it the active sheet/tab
$ws.Range("I3").Select()
$ws.Activate()
